$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert two new rows for the new "lag1" / "lag1+PCA" model results ---
$null = $ws.Rows("6:7").Insert(-4121)

# --- Row 6: Lasso Regression+normalization+lag1 ---
$ws.Range("A6").Value = 4
$ws.Range("B6").Value = "Lasso Regression+normalization+lag1"
$ws.Range("C6").Value = 81.344944740947696
$ws.Range("D6").Value = 97.302968419999999
$ws.Range("E6").Value = 69.779661556946493
$ws.Range("F6").Value = 77.074776640014207
$ws.Range("G6").Value = 71.002143129999993

# --- Row 7: Lasso Regression+normalization+lag1+PCA ---
$ws.Range("A7").Value = 5
$ws.Range("B7").Value = "Lasso Regression+normalization+lag1+PCA"
$ws.Range("C7").Value = 88.856083384300604
$ws.Range("D7").Value = 92.061861870000001
$ws.Range("E7").Value = 73.409172663976904
$ws.Range("F7").Value = 85.044508627085506
$ws.Range("G7").Value = 81.644887359999998

# --- Give the two new rows the same "boxed" medium border used elsewhere in the table ---
foreach ($r in 6, 7) {
    $ws.Range("A$r").Borders.Item(7).Weight = -4138
    $ws.Range("A$r").Borders.Item(10).Weight = -4138
    $ws.Range("B$r").Borders.Item(7).Weight = -4138
    $ws.Range("C$r").Borders.Item(7).Weight = -4138
    $ws.Range("C$r").Borders.Item(10).Weight = -4138
    $ws.Range("E$r").Borders.Item(7).Weight = -4138
    $ws.Range("E$r").Borders.Item(10).Weight = -4138
    $ws.Range("F$r").Borders.Item(7).Weight = -4138
    $ws.Range("F$r").Borders.Item(10).Weight = -4138
}
$ws.Range("D6:G7").Font.Name = "DengXian"

# --- Renumber the Id column for the rows that shifted down two places ---
$ws.Range("A8").Value = 6
$ws.Range("A9").Value = 7
$ws.Range("A10").Value = 8
$ws.Range("A11").Value = 9
$ws.Range("A12").Value = 10

# --- Widen column B to fit the longer model names ---
$ws.Columns("B").ColumnWidth = 39.6

# --- Restore the on-screen selection ---
$null = $ws.Range("F16").Select()
